$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "level2" / "level3" / "level4" sheets by copying the
#    existing "level1" sheet (this keeps formatting / styles identical,
#    exactly like the author did when building the new tabs).
# ---------------------------------------------------------------------------
$level1 = $wb.Worksheets.Item("level1")

$level1.Copy([System.Reflection.Missing]::Value, $level1)
$level2 = $wb.Worksheets.Item($level1.Index + 1)
$level2.Name = "level2"

$level2.Copy([System.Reflection.Missing]::Value, $level2)
$level3 = $wb.Worksheets.Item($level2.Index + 1)
$level3.Name = "level3"

$level3.Copy([System.Reflection.Missing]::Value, $level3)
$level4 = $wb.Worksheets.Item($level3.Index + 1)
$level4.Name = "level4"

# ---------------------------------------------------------------------------
# 2. access
# ---------------------------------------------------------------------------
$access = $wb.Worksheets.Item("access")
$access.Range("B2").Value = " Analyse des enfants accédant à l'éducation"
$access.Range("B10").Select()

# ---------------------------------------------------------------------------
# 3. out_of_school
# ---------------------------------------------------------------------------
$oos = $wb.Worksheets.Item("out_of_school")
$oos.Range("B2").Value = "Analyse des enfants n'ayant pas accès à l'éducation, OoS"
$oos.Range("G2").Value = "child-marriage"
$ps = $oos.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
$oos.Range("B10").Select()

# ---------------------------------------------------------------------------
# 4. ece
# ---------------------------------------------------------------------------
$ece = $wb.Worksheets.Item("ece")
$ece.Range("B2").Value = "éducation pré-primaire"
$ece.Range("G2").Value = "% d'enfants un an avant l'âge d'entrée à l'école primaire fréquentant un programme d'éducation préscolaire ou l'école primaire"
$ece.Range("G3").Value = "% d'enfants un an avant l'âge d'entrée à l'école primaire fréquentant l'école primaire"
$ece.Range("D9").Select()

# ---------------------------------------------------------------------------
# 5. level1 (primaire)
# ---------------------------------------------------------------------------
$level1.Range("B2").Value = "Profil de fréquentation scolaire: primaire"
$level1.Range("E10").Select()

# ---------------------------------------------------------------------------
# 6. level2 (niveau intermédiaire -- secondaire)
# ---------------------------------------------------------------------------
$level2.Range("B2").Value = "Profil de fréquentation scolaire: niveau intermédiaire -- secondaire"
$level2.Range("G2").Value = "% d'enfants en âge de niveau intermédiaire actuellement scolarisés au niveau intermédiaire ou plus"
$level2.Range("G3").Value = "% d'enfants avec 2 ans de plus que l'âge prévu : secondaire"
$level2.Range("B2").Select()

# ---------------------------------------------------------------------------
# 7. level3 (secondaire)
# ---------------------------------------------------------------------------
$level3.Range("B2").Value = "Profil de fréquentation scolaire: secondaire"
$level3.Range("G2").Value = "% d'enfants en âge de secondaire actuellement scolarisés au secondaire ou plus"
$level3.Range("G3").ClearContents()
$level3.Range("E7").Select()

# ---------------------------------------------------------------------------
# 8. level4 (secondaire supérieur)
# ---------------------------------------------------------------------------
$level4.Range("B2").Value = "Profil de fréquentation scolaire:  secondaire supérieur"
$level4.Range("G2").Value = "% d'enfants en âge de secondaire supérieur actuellement scolarisés au secondaire supérieur ou plus"
$level4.Range("G3").ClearContents()
$level4.Range("E9").Select()

# ---------------------------------------------------------------------------
# 9. Restore "out_of_school" as the active/visible tab.
# ---------------------------------------------------------------------------
$oos.Activate()
